$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price (D) column as Text so strings containing dots/
# percent-like patterns are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.928.04"

$ws.Range("D3").Value = "1.812.77"
$ws.Range("E3").Value = "  +1.93%  "

$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("D5").Value = "311.69"
$ws.Range("E5").Value = "  +1.42%  "

$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("D7").Value = "0.4286"
$ws.Range("E7").Value = "  +1.47%  "

$ws.Range("D8").Value = "0.3685"
$ws.Range("E8").Value = "  +1.89%  "

$ws.Range("D9").Value = "0.07241"
$ws.Range("E9").Value = "  +1.00%  "

$ws.Range("D10").Value = "0.8612"
$ws.Range("E10").Value = "  +2.79%  "

$ws.Range("D11").Value = "2.050.77"
$ws.Range("E11").Value = "  +16.66%  "

$ws.Range("D12").Value = "21.16"
$ws.Range("E12").Value = "  +4.45%  "

$ws.Range("D13").Value = "6.630"
$ws.Range("E13").Value = "  +4.57%  "

$ws.Range("D14").Value = "5.381"
$ws.Range("E14").Value = "  +2.49%  "

$ws.Range("D15").Value = "0.06890"
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("E16").Value = "  +1.90%  "

$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  -0.70%  "

$ws.Range("D18").Value = "0.000008855"
$ws.Range("E18").Value = "  +2.02%  "

$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("E20").Value = "  +1.58%  "

$ws.Range("D21").Value = "26.965.35"
$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("D22").Value = "5.190"
$ws.Range("E22").Value = "  +3.57%  "

$ws.Range("D23").Value = "11.02"
$ws.Range("E23").Value = "  -0.47%  "

$ws.Range("D24").Value = "2.286.88"
$ws.Range("E24").Value = "  +15.43%  "

$ws.Range("D25").Value = "153.79"
$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("D26").Value = "1.885"
$ws.Range("E26").Value = "  -1.90%  "

$ws.Range("D27").Value = "18.30"
$ws.Range("E27").Value = "  +0.87%  "

$ws.Range("D28").Value = "5.211"
$ws.Range("E28").Value = "  +3.32%  "

$ws.Range("D29").Value = "1.890"
$ws.Range("E29").Value = "  +16.30%  "

$ws.Range("E30").Value = "  +0.47%  "

$ws.Range("D31").Value = "0.08931"
$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("D32").Value = "0.7426"
$ws.Range("E32").Value = "  +3.17%  "

$ws.Range("D33").Value = "1.160"
$ws.Range("E33").Value = "  +6.40%  "

$ws.Range("E34").Value = "  +2.36%  "

$ws.Range("D35").Value = "2.795"
$ws.Range("E35").Value = "  -1.47%  "

$ws.Range("D37").Value = "1.122"
$ws.Range("E37").Value = "  +3.83%  "

$ws.Range("D38").Value = "0.05212"
$ws.Range("E38").Value = "  +2.57%  "

$ws.Range("D39").Value = "0.01925"
$ws.Range("E39").Value = "  +1.74%  "

$ws.Range("D40").Value = "0.5086"
$ws.Range("E40").Value = "  +3.45%  "

$ws.Range("D41").Value = "2.758"
$ws.Range("E41").Value = "  +9.50%  "

$ws.Range("D42").Value = "0.1640"
$ws.Range("E42").Value = "  +1.92%  "

$ws.Range("D43").Value = "6.436"
$ws.Range("E43").Value = "  +5.85%  "

$ws.Range("D44").Value = "8.260"
$ws.Range("E44").Value = "  +4.28%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "106.79"
$ws.Range("E45").Value = "  +2.25%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "10.40"
$ws.Range("E46").Value = "  +3.23%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "1.003"
$ws.Range("E47").Value = "  -0.43%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "0.4582"
$ws.Range("E48").Value = "  +2.34%  "

$ws.Range("D49").Value = "1.651"
$ws.Range("E49").Value = "  +5.12%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.06286"
$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "1.806"
$ws.Range("E51").Value = "  +5.07%  "

# Restore default (General) formatting on the Price column so the
# underlying cell style matches the original (no explicit number format),
# while keeping the values as text strings.
$ws.Range("D2:D51").ClearFormats()
